# Ptprz1-Ncam1.xlsx was regenerated with new TPM (transcripts-per-million)
# input data. The ligand/receptor pair (Ptprz1/Ncam1) is unchanged, but the
# "Sending cluster" x "Target cluster" grid grew from 2 senders x 3 targets
# (rows 2-7) to 3 senders x 3 targets (rows 2-10: ECs/FAPs/MuSCs each paired
# with ECs/FAPs/MuSCs), and every numeric column was recomputed.
#
# Rewrite the full data block (rows 2-10, columns A:T) with the values from
# the updated pipeline run. Column order per row:
#   A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
#   E..J = ligand-expressing-cell stats, K..P = receptor-expressing-cell
#   stats, Q..T = edge weight / specificity stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs",   "Ptprz1", "Ncam1", "ECs",   2, 0.6666666666666666, 0.06163433333333333, 0.184903,           0.006690894379667537, 0.006690894379667537, 2, 0.6666666666666666, 0.04741066666666666, 0.142232,          0.003188134523263584, 0.003188134523263585, 0.002922124832888888, 0.026299123496,    0.00002133147136332836, 0.00002133147136332836),
  @("ECs",   "Ptprz1", "Ncam1", "FAPs",  2, 0.6666666666666666, 0.06163433333333333, 0.184903,           0.006690894379667537, 0.006690894379667537, 3, 1,                  0.237305,            0.711915,          0.01595759596384214,  0.01595759596384214,  0.01462613547166666,  0.131635219245,    0.0001067705891474767,  0.0001067705891474767),
  @("ECs",   "Ptprz1", "Ncam1", "MuSCs", 2, 0.6666666666666666, 0.06163433333333333, 0.184903,           0.006690894379667537, 0.006690894379667537, 3, 1,                  14.58625866666667,   43.758776,         0.9808542695128942,  0.9808542695128943,  0.8990143287475554,  8.091128958728,    0.006562792319156732,   0.006562792319156733),
  @("FAPs",  "Ptprz1", "Ncam1", "ECs",   3, 1,                  0.1030763333333333, 0.309229,           0.01118975126488057, 0.01118975126488057,  2, 0.6666666666666666, 0.04741066666666666, 0.142232,          0.003188134523263584, 0.003188134523263585, 0.004886917680888888, 0.043982259128,    0.00003567443231429812, 0.00003567443231429812),
  @("FAPs",  "Ptprz1", "Ncam1", "FAPs",  3, 1,                  0.1030763333333333, 0.309229,           0.01118975126488057, 0.01118975126488057,  3, 1,                  0.237305,            0.711915,          0.01595759596384214,  0.01595759596384214,  0.02446052928166666,  0.220144763535,    0.0001785615296208557,  0.0001785615296208557),
  @("FAPs",  "Ptprz1", "Ncam1", "MuSCs", 3, 1,                  0.1030763333333333, 0.309229,           0.01118975126488057, 0.01118975126488057,  3, 1,                  14.58625866666667,   43.758776,         0.9808542695128942,  0.9808542695128943,  1.503498060411555,   13.531482543704,   0.01097551530294542,    0.01097551530294542),
  @("MuSCs", "Ptprz1", "Ncam1", "ECs",   3, 1,                  9.046962666666667,  27.140888,          0.982119354355452,   0.9821193543554518,   2, 0.6666666666666666, 0.04741066666666666, 0.142232,          0.003188134523263584, 0.003188134523263585, 0.4289225313351111,  3.860302782016,    0.003131128619585958,   0.003131128619585958),
  @("MuSCs", "Ptprz1", "Ncam1", "FAPs",  3, 1,                  9.046962666666667,  27.140888,          0.982119354355452,   0.9821193543554518,   3, 1,                  0.237305,            0.711915,          0.01595759596384214,  0.01595759596384214,  2.146889475613333,   19.32200528052,    0.0156722638450738,     0.0156722638450738),
  @("MuSCs", "Ptprz1", "Ncam1", "MuSCs", 3, 1,                  9.046962666666667,  27.140888,          0.982119354355452,   0.9821193543554518,   3, 1,                  14.58625866666667,   43.758776,         0.9808542695128942,  0.9808542695128943,  131.9613376036764,   1187.652038433088, 0.9633159618907922,     0.9633159618907922)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $rowvals = $data[$i]
  for ($j = 0; $j -lt $rowvals.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value2 = $rowvals[$j]
  }
}
